# Weekly fruit/vegetable price update: insert a new weekly record
# (Vega Monumental Concepción - Piña - Caramelo - Primera) at row 152,
# pushing the existing rows 152-240 down to 153-241.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 152 (shifts rows 152..240 -> 153..241)
$ws.Rows.Item(152).Insert()

# Populate the new row 152 with this week's data point
$ws.Range("A152").Value2 = 11
$ws.Range("B152").Value2 = "Vega Monumental Concepción"
$ws.Range("C152").Value2 = "Bíobío"
$ws.Range("D152").Value2 = 44960
$ws.Range("E152").Value2 = 8
$ws.Range("F152").Value2 = "Fruta"
$ws.Range("G152").Value2 = 100108
$ws.Range("H152").Value2 = "Tropicales y subtropicales"
$ws.Range("I152").Value2 = 100108005
$ws.Range("J152").Value2 = "Piña"
$ws.Range("K152").Value2 = "Caramelo"
$ws.Range("L152").Value2 = "Primera"
$ws.Range("M152").Value2 = 200
$ws.Range("N152").Value2 = 19000
$ws.Range("O152").Value2 = 20000
$ws.Range("P152").Value2 = 19500
$ws.Range("Q152").Value2 = "$/caja 12 unidades"
$ws.Range("R152").Value2 = "Ecuador"
$ws.Range("S152").Value2 = 1625
$ws.Range("T152").Value2 = 12
